$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.919.94"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "1.812.72"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.24"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4644"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3653"
$ws.Range("E8").Value = "  -1.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07358"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8679"
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.24"
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("D12").Value = "1.814.05"
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.359"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07092"
$ws.Range("E14").Value = "  +0.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.495"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.08"
$ws.Range("E16").Value = "  -1.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008706"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.63"
$ws.Range("D21").Value = "26.936.38"
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.282"
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("E23").Value = "  -0.47%  "
$ws.Range("D24").Value = "2.051.20"
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.896"
$ws.Range("E25").Value = "  -0.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "150.77"
$ws.Range("E26").Value = "  -0.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.31"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.121"
$ws.Range("E28").Value = "  -1.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.245"
$ws.Range("E29").Value = "  -0.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.46"
$ws.Range("E30").Value = "  -0.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08877"
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7537"
$ws.Range("E32").Value = "  -0.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.163"
$ws.Range("E33").Value = "  +0.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.473"
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.904"
$ws.Range("E35").Value = "  -0.83%  "
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("E37").Value = "  -1.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05283"
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01944"
$ws.Range("E39").Value = "  -1.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.977"
$ws.Range("E40").Value = "  +1.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.249"
$ws.Range("E41").Value = "  +0.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5297"
$ws.Range("E42").Value = "  -0.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.289"
$ws.Range("E43").Value = "  -4.76%  "
$ws.Range("E44").Value = "  -0.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.401"
$ws.Range("E45").Value = "  -1.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4860"
$ws.Range("E46").Value = "  -2.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.44"
$ws.Range("E47").Value = "  +0.21%  "
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.660"
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.01"
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("E51").Value = "  -0.08%  "
